# Update the worksheet name to reflect the new extraction timestamp
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Name = "IClientBalance-20240528-092752-"

# Update the date column (G2:G257) from 45439 (2024-05-27) to 45440 (2024-05-28)
$ws.Range("G2:G257").Value = 45440
